$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 62: set H62=0, I62=0, J62=0, K62=0, L62=0
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
# row 62: clear M62:N62
$ws.Range("M62:N62").ClearContents()
# row 65: set H65=0, I65=0, J65=0, K65=0, L65=0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
# row 65: clear M65:N65
$ws.Range("M65:N65").ClearContents()
# row 86: set H86=55559830, I86=71431930, J86=7500, K86=71431930, L86=7500, M86=-71430807, N86=-9746
$ws.Range("H86").Value = 55559830
$ws.Range("I86").Value = 71431930
$ws.Range("J86").Value = 7500
$ws.Range("K86").Value = 71431930
$ws.Range("L86").Value = 7500
$ws.Range("M86").Value = -71430807
$ws.Range("N86").Value = -9746
# row 89: set H89=55559830, I89=71431930, J89=7500, K89=357159650, L89=37500, M89=-357154034, N89=-48732
$ws.Range("H89").Value = 55559830
$ws.Range("I89").Value = 71431930
$ws.Range("J89").Value = 7500
$ws.Range("K89").Value = 357159650
$ws.Range("L89").Value = 37500
$ws.Range("M89").Value = -357154034
$ws.Range("N89").Value = -48732
# row 100: set H100=1810.8096, I100=1470.9, J100=2119.818, K100=1470.9, L100=2119.818, M100=-929.9000000000001, N100=-3201.818
$ws.Range("H100").Value = 1810.8096
$ws.Range("I100").Value = 1470.9
$ws.Range("J100").Value = 2119.818
$ws.Range("K100").Value = 1470.9
$ws.Range("L100").Value = 2119.818
$ws.Range("M100").Value = -929.9000000000001
$ws.Range("N100").Value = -3201.818
# row 113: set H113=2891.4062, I113=2366.25, J113=3416.5625, K113=2366.25, L113=3416.5625, M113=887.75, N113=-9924.5625
$ws.Range("H113").Value = 2891.4062
$ws.Range("I113").Value = 2366.25
$ws.Range("J113").Value = 3416.5625
$ws.Range("K113").Value = 2366.25
$ws.Range("L113").Value = 3416.5625
$ws.Range("M113").Value = 887.75
$ws.Range("N113").Value = -9924.5625
# row 127: set H127=1503.871, I127=467.2857, K127=1401.8571, M127=3558.1429
$ws.Range("H127").Value = 1503.871
$ws.Range("I127").Value = 467.2857
$ws.Range("K127").Value = 1401.8571
$ws.Range("M127").Value = 3558.1429
# row 137: set H137=13048.739, I137=1007.6667, K137=3023.0001, M137=-473.0001000000002
$ws.Range("H137").Value = 13048.739
$ws.Range("I137").Value = 1007.6667
$ws.Range("K137").Value = 3023.0001
$ws.Range("M137").Value = -473.0001000000002
# row 138: set H138=7411574, I138=14495477, J138=5675, K138=43486431, L138=17025, M138=-43481291, N138=-27305
$ws.Range("H138").Value = 7411574
$ws.Range("I138").Value = 14495477
$ws.Range("J138").Value = 5675
$ws.Range("K138").Value = 43486431
$ws.Range("L138").Value = 17025
$ws.Range("M138").Value = -43481291
$ws.Range("N138").Value = -27305

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 32: set H32=18729.475, I32=18729.475, K32=18729.475, M32=-18442.475
$ws.Range("H32").Value = 18729.475
$ws.Range("I32").Value = 18729.475
$ws.Range("K32").Value = 18729.475
$ws.Range("M32").Value = -18442.475
# row 61: set H61=1814.9824, I61=1758.1072, K61=1758.1072, M61=-1546.1072
$ws.Range("H61").Value = 1814.9824
$ws.Range("I61").Value = 1758.1072
$ws.Range("K61").Value = 1758.1072
$ws.Range("M61").Value = -1546.1072
# row 110: set H110=825, I110=433.33334, J110=2000, K110=433.33334, L110=2000, M110=1611.66666, N110=-6090
$ws.Range("H110").Value = 825
$ws.Range("I110").Value = 433.33334
$ws.Range("J110").Value = 2000
$ws.Range("K110").Value = 433.33334
$ws.Range("L110").Value = 2000
$ws.Range("M110").Value = 1611.66666
$ws.Range("N110").Value = -6090
# row 136: set H136=1814.9824, I136=1758.1072, K136=5274.321599999999, M136=-2724.321599999999
$ws.Range("H136").Value = 1814.9824
$ws.Range("I136").Value = 1758.1072
$ws.Range("K136").Value = 5274.321599999999
$ws.Range("M136").Value = -2724.321599999999

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 31: set H31=6375.306, I31=2346.7058, J31=8515.5, K31=2346.7058, L31=8515.5, M31=-2051.7058, N31=-9105.5
$ws.Range("H31").Value = 6375.306
$ws.Range("I31").Value = 2346.7058
$ws.Range("J31").Value = 8515.5
$ws.Range("K31").Value = 2346.7058
$ws.Range("L31").Value = 8515.5
$ws.Range("M31").Value = -2051.7058
$ws.Range("N31").Value = -9105.5
# row 34: set H34=6375.306, I34=2346.7058, J34=8515.5, K34=2346.7058, L34=8515.5, M34=-2144.7058, N34=-8919.5
$ws.Range("H34").Value = 6375.306
$ws.Range("I34").Value = 2346.7058
$ws.Range("J34").Value = 8515.5
$ws.Range("K34").Value = 2346.7058
$ws.Range("L34").Value = 8515.5
$ws.Range("M34").Value = -2144.7058
$ws.Range("N34").Value = -8919.5
# row 48: set H48=23760.4, J48=23760.4, L48=23760.4, N48=-24712.4
$ws.Range("H48").Value = 23760.4
$ws.Range("J48").Value = 23760.4
$ws.Range("L48").Value = 23760.4
$ws.Range("N48").Value = -24712.4
# row 58: set H58=1736.3529, I58=722.6957, J58=3855.818, K58=722.6957, L58=3855.818, M58=-519.6957, N58=-4261.818
$ws.Range("H58").Value = 1736.3529
$ws.Range("I58").Value = 722.6957
$ws.Range("J58").Value = 3855.818
$ws.Range("K58").Value = 722.6957
$ws.Range("L58").Value = 3855.818
$ws.Range("M58").Value = -519.6957
$ws.Range("N58").Value = -4261.818
# row 99: set H99=3149.5, I99=1575, J99=3936.75, K99=1575, L99=3936.75, M99=-77, N99=-6932.75
$ws.Range("H99").Value = 3149.5
$ws.Range("I99").Value = 1575
$ws.Range("J99").Value = 3936.75
$ws.Range("K99").Value = 1575
$ws.Range("L99").Value = 3936.75
$ws.Range("M99").Value = -77
$ws.Range("N99").Value = -6932.75
# row 107: set H107=564.7778, I107=369.78262, J107=1686, K107=369.78262, L107=1686, M107=1550.21738, N107=-5526
$ws.Range("H107").Value = 564.7778
$ws.Range("I107").Value = 369.78262
$ws.Range("J107").Value = 1686
$ws.Range("K107").Value = 369.78262
$ws.Range("L107").Value = 1686
$ws.Range("M107").Value = 1550.21738
$ws.Range("N107").Value = -5526
# row 122: set H122=1738.2632, I122=1886.4615, J122=1417.1666, K122=5659.3845, L122=4251.4998, M122=-3209.3845, N122=-9151.4998
$ws.Range("H122").Value = 1738.2632
$ws.Range("I122").Value = 1886.4615
$ws.Range("J122").Value = 1417.1666
$ws.Range("K122").Value = 5659.3845
$ws.Range("L122").Value = 4251.4998
$ws.Range("M122").Value = -3209.3845
$ws.Range("N122").Value = -9151.4998
# row 126: set H126=3149.5, I126=1575, J126=3936.75, K126=4725, L126=11810.25, M126=-2255, N126=-16750.25
$ws.Range("H126").Value = 3149.5
$ws.Range("I126").Value = 1575
$ws.Range("J126").Value = 3936.75
$ws.Range("K126").Value = 4725
$ws.Range("L126").Value = 11810.25
$ws.Range("M126").Value = -2255
$ws.Range("N126").Value = -16750.25
# row 136: set H136=1736.3529, I136=722.6957, J136=3855.818, K136=2168.0871, L136=11567.454, M136=381.9129000000003, N136=-16667.454
$ws.Range("H136").Value = 1736.3529
$ws.Range("I136").Value = 722.6957
$ws.Range("J136").Value = 3855.818
$ws.Range("K136").Value = 2168.0871
$ws.Range("L136").Value = 11567.454
$ws.Range("M136").Value = 381.9129000000003
$ws.Range("N136").Value = -16667.454

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 37: set H37=49900, J37=49900, L37=149700, N37=-149924
$ws.Range("H37").Value = 49900
$ws.Range("J37").Value = 49900
$ws.Range("L37").Value = 149700
$ws.Range("N37").Value = -149924
# row 133: set H133=55560308, I133=100002230, J133=7900, K133=300006690, L133=23700, M133=-300001630, N133=-33820
$ws.Range("H133").Value = 55560308
$ws.Range("I133").Value = 100002230
$ws.Range("J133").Value = 7900
$ws.Range("K133").Value = 300006690
$ws.Range("L133").Value = 23700
$ws.Range("M133").Value = -300001630
$ws.Range("N133").Value = -33820

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 102: set H102=3662.5652, I102=5495.6, J102=2252.5386, K102=5495.6, L102=2252.5386, M102=-3873.6, N102=-5496.5386
$ws.Range("H102").Value = 3662.5652
$ws.Range("I102").Value = 5495.6
$ws.Range("J102").Value = 2252.5386
$ws.Range("K102").Value = 5495.6
$ws.Range("L102").Value = 2252.5386
$ws.Range("M102").Value = -3873.6
$ws.Range("N102").Value = -5496.5386
# row 126: set H126=4883.0835, I126=0, J126=4883.0835, K126=0, L126=14649.2505, N126=-19589.2505
$ws.Range("H126").Value = 4883.0835
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 4883.0835
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 14649.2505
$ws.Range("N126").Value = -19589.2505
# row 126: clear M126
$ws.Range("M126").ClearContents()

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 16: set H16=2968.0715, I16=3187.625, J16=2675.3333, K16=3187.625, L16=2675.3333, M16=-3017.625, N16=-3015.3333
$ws.Range("H16").Value = 2968.0715
$ws.Range("I16").Value = 3187.625
$ws.Range("J16").Value = 2675.3333
$ws.Range("K16").Value = 3187.625
$ws.Range("L16").Value = 2675.3333
$ws.Range("M16").Value = -3017.625
$ws.Range("N16").Value = -3015.3333
# row 55: set H55=318.29166, I55=170.84616, J55=492.54544, K55=170.84616, L55=492.54544, M55=2.153840000000002, N55=-838.54544
$ws.Range("H55").Value = 318.29166
$ws.Range("I55").Value = 170.84616
$ws.Range("J55").Value = 492.54544
$ws.Range("K55").Value = 170.84616
$ws.Range("L55").Value = 492.54544
$ws.Range("M55").Value = 2.153840000000002
$ws.Range("N55").Value = -838.54544
# row 68: set H68=2160.5217, I68=1952.5294, K68=1952.5294, M68=-1203.5294
$ws.Range("H68").Value = 2160.5217
$ws.Range("I68").Value = 1952.5294
$ws.Range("K68").Value = 1952.5294
$ws.Range("M68").Value = -1203.5294
# row 71: set H71=2160.5217, I71=1952.5294, K71=9762.646999999999, M71=-6018.646999999999
$ws.Range("H71").Value = 2160.5217
$ws.Range("I71").Value = 1952.5294
$ws.Range("K71").Value = 9762.646999999999
$ws.Range("M71").Value = -6018.646999999999
# row 132: set H132=10423932, I132=5631.125, J132=31260534, K132=16893.375, L132=93781602, M132=-14363.375, N132=-93786662
$ws.Range("H132").Value = 10423932
$ws.Range("I132").Value = 5631.125
$ws.Range("J132").Value = 31260534
$ws.Range("K132").Value = 16893.375
$ws.Range("L132").Value = 93781602
$ws.Range("M132").Value = -14363.375
$ws.Range("N132").Value = -93786662
# row 136: set H136=4997.683, I136=2293.7097, J136=13380, K136=6881.1291, L136=40140, M136=-4331.1291, N136=-45240
$ws.Range("H136").Value = 4997.683
$ws.Range("I136").Value = 2293.7097
$ws.Range("J136").Value = 13380
$ws.Range("K136").Value = 6881.1291
$ws.Range("L136").Value = 40140
$ws.Range("M136").Value = -4331.1291
$ws.Range("N136").Value = -45240

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 62: set H62=5937.4, I62=5755.5557, J62=6086.1816, K62=5755.5557, L62=6086.1816, M62=-5131.5557, N62=-7334.1816
$ws.Range("H62").Value = 5937.4
$ws.Range("I62").Value = 5755.5557
$ws.Range("J62").Value = 6086.1816
$ws.Range("K62").Value = 5755.5557
$ws.Range("L62").Value = 6086.1816
$ws.Range("M62").Value = -5131.5557
$ws.Range("N62").Value = -7334.1816
# row 65: set H65=5937.4, I65=5755.5557, J65=6086.1816, K65=28777.7785, L65=30430.908, M65=-25657.7785, N65=-36670.908
$ws.Range("H65").Value = 5937.4
$ws.Range("I65").Value = 5755.5557
$ws.Range("J65").Value = 6086.1816
$ws.Range("K65").Value = 28777.7785
$ws.Range("L65").Value = 30430.908
$ws.Range("M65").Value = -25657.7785
$ws.Range("N65").Value = -36670.908
# row 81: set H81=2066.6667, I81=500, K81=1000, M81=61
$ws.Range("H81").Value = 2066.6667
$ws.Range("I81").Value = 500
$ws.Range("K81").Value = 1000
$ws.Range("M81").Value = 61
# row 84: set H84=2066.6667, I84=500, K84=5000, M84=304
$ws.Range("H84").Value = 2066.6667
$ws.Range("I84").Value = 500
$ws.Range("K84").Value = 5000
$ws.Range("M84").Value = 304
# row 122: set H122=2066.6667, I122=2530.2778, J122=1139.4445, K122=7590.8334, L122=3418.3335, M122=-5140.8334, N122=-8318.333500000001
$ws.Range("H122").Value = 2066.6667
$ws.Range("I122").Value = 2530.2778
$ws.Range("J122").Value = 1139.4445
$ws.Range("K122").Value = 7590.8334
$ws.Range("L122").Value = 3418.3335
$ws.Range("M122").Value = -5140.8334
$ws.Range("N122").Value = -8318.333500000001
# row 126: set H126=2026.921, I126=2021.7368, J126=2032.1052, K126=6065.2104, L126=6096.3156, M126=-3595.2104, N126=-11036.3156
$ws.Range("H126").Value = 2026.921
$ws.Range("I126").Value = 2021.7368
$ws.Range("J126").Value = 2032.1052
$ws.Range("K126").Value = 6065.2104
$ws.Range("L126").Value = 6096.3156
$ws.Range("M126").Value = -3595.2104
$ws.Range("N126").Value = -11036.3156
# row 132: set H132=3410.9092, I132=2841.4443, J132=5973.5, K132=8524.332900000001, L132=17920.5, M132=-5994.332900000001, N132=-22980.5
$ws.Range("H132").Value = 3410.9092
$ws.Range("I132").Value = 2841.4443
$ws.Range("J132").Value = 5973.5
$ws.Range("K132").Value = 8524.332900000001
$ws.Range("L132").Value = 17920.5
$ws.Range("M132").Value = -5994.332900000001
$ws.Range("N132").Value = -22980.5
